$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览" (exhibitions) - update column F (想去人数 / "want to go" count)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 939
$ws1.Range("F6").Value  = 4879
$ws1.Range("F7").Value  = 376
$ws1.Range("F8").Value  = 554
$ws1.Range("F9").Value  = 858
$ws1.Range("F16").Value = 1576
$ws1.Range("F18").Value = 671
$ws1.Range("F23").Value = 108
$ws1.Range("F27").Value = 1405
$ws1.Range("F29").Value = 73
$ws1.Range("F30").Value = 10
$ws1.Range("F37").Value = 535
$ws1.Range("F40").Value = 11

# Sheet 4 = "全部类型" (all types) - same underlying rows, mirrored layout
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value  = 939
$ws4.Range("F8").Value  = 4879
$ws4.Range("F9").Value  = 376
$ws4.Range("F10").Value = 554
$ws4.Range("F13").Value = 858
$ws4.Range("F23").Value = 1576
$ws4.Range("F25").Value = 671
$ws4.Range("F31").Value = 108
$ws4.Range("F34").Value = 1405
$ws4.Range("F36").Value = 73
$ws4.Range("F37").Value = 10
$ws4.Range("F43").Value = 535
$ws4.Range("F46").Value = 11
